# Continue the paragraph that ends with "... ze stejného majetku." by
# adding two new sentences in place of the stray "_GoBack" bookmark that
# Word leaves at the last edit position, then remove that bookmark.

$d = $word.ActiveDocument

$bm = $null
foreach ($b in $d.Bookmarks) {
    if ($b.Name -eq "_GoBack") {
        $bm = $b
    }
}
if (-not $bm) {
    $bm = $d.Bookmarks.Item("_GoBack")
}

$r = $bm.Range

# First new sentence fragment (keeps the leading space like the original run).
$r.InsertAfter(" Může například chtít")
$r.Collapse(0)

# Second new sentence, appended right after the first.
$r.InsertAfter(", aby zakladatelské právní jednání bylo jednostranné.")
$r.Collapse(0)

# The bookmark has served its purpose (matching the diff, which drops it).
$bm2 = $d.Bookmarks.Item("_GoBack")
$bm2.Delete()
